$wb = $excel.ActiveWorkbook

# --- Update the dataset_type lookup sheet to accommodate the newly introduced
#     "2D Imaging Mass Cytometry" assay: drop the retired "nanoPOTS" and
#     "NanoDESI" assays, and add the new assay row right after "MALDI".
$ds = $wb.Worksheets.Item("dataset_type")

# Row 3 = nanoPOTS -> remove entirely (whole row shifts up)
$ds.Rows.Item(3).Delete()

# After the delete above, NanoDESI (was row 21) is now row 20 -> remove it too
$ds.Rows.Item(20).Delete()

# After both deletes, MALDI (was row 24) is now row 22, and RNAseq (GeoMx)
# (was row 25) is now row 23. Insert a new blank row before RNAseq (GeoMx)
# i.e. right after MALDI, and populate it with the new assay.
$ds.Rows.Item(23).Insert()
$ds.Cells.Item(23, 1).Value = "2D Imaging Mass Cytometry"
$ds.Cells.Item(23, 2).Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296"

# The dataset_type list now has 35 entries (36 - nanoPOTS - NanoDESI + new
# assay) instead of 36, so the dropdown validation on the main sheet needs to
# point at the smaller range.
$main = $wb.Worksheets.Item("Visium")
$main.Range("D2:D1001").Validation.Formula1 = "'dataset_type'!`$A`$1:`$A`$35"

# --- Bump the metadata "createdOn" timestamp on the .metadata sheet.
$meta = $wb.Worksheets.Item(".metadata")
$meta.Cells.Item(2, 3).Value = "2023-11-02T15:46:14-07:00"
